$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the existing "Certificado" / certificate-number pair by moving it
# (value + style, incl. the Text number format) from column A into column C,
# where it now belongs as part of the wider intake form.
$ws.Range("A1:A2").Copy($ws.Range("C1:C2"))

# Column A no longer carries the old "Text" style for the certificate number;
# it becomes the plain "User" credential cell.
$ws.Range("A2").ClearFormats()

# Fill in the rest of the intake form, in the same left-to-right / row-by-row
# order the form was originally completed in.
$ws.Range("A1").Value = "User"
$ws.Range("B1").Value = "Password"
$ws.Range("A2").Value = "metlifecustomeruser"
$ws.Range("B2").Value = "MetLifePilot@1234"

$ws.Range("D1").Value = "TypeOfClaim"
$ws.Range("D2").Value = "Complementary"

$ws.Range("E1").Value = "Abroad"
$ws.Range("E2").Value = "Yes"

$ws.Range("F1").Value = "MetEmployee"
$ws.Range("F2").Value = "No"

$ws.Range("G2").Value = "Hombro Roto"
$ws.Range("G1").Value = "IllnessDetails"

$ws.Range("H1").Value = "Email"
$ws.Range("I1").Value = "ContactNumber"
$ws.Range("H2").Value = "isidrocarrasco@prueba.com"

$ws.Range("J1").Value = "Hospital"
$ws.Range("J2").Value = "HOSPITAL ANGELES METROPOLITANO"

$ws.Range("K2").Value = "Vianney Roman Garcia"
$ws.Range("K1").Value = "Doctor"

$ws.Range("L1").Value = "Adicional"
$ws.Range("L2").Value = "Laboratorio"

# Numeric contact number (kept as a real number, not text).
$ws.Range("I2").Value = 6730776222

# Column widths to accommodate the new data.
$ws.Columns.Item(1).ColumnWidth = 18.666667
$ws.Columns.Item(2).ColumnWidth = 17.333333
$ws.Columns.Item(3).ColumnWidth = 14.166667
$ws.Columns.Item(4).ColumnWidth = 14.666667
$ws.Columns.Item(6).ColumnWidth = 12.833333
$ws.Columns.Item(7).ColumnWidth = 12.666667
$ws.Columns.Item(8).ColumnWidth = 25.166667
$ws.Columns.Item(9).ColumnWidth = 14.666667
$ws.Columns.Item(10).ColumnWidth = 32.666667
$ws.Columns.Item(11).ColumnWidth = 21.666667

# The certificate number is numeric-looking text; tell Excel to ignore the
# "number stored as text" warning for it, like the author did.
try {
    $ws.Range("C2").Errors.Item(9).Ignore = $true
} catch {
}

# Restore the view: scrolled so column F is leftmost, with L2 selected.
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$win.ScrollRow = 1
$ws.Range("L2").Select()
